$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '43.735.04'
$ws.Range("E2").Value = '  +2.15%  '

$ws.Range("D3").Value = '2.215.15'
$ws.Range("E3").Value = '  +0.08%  '

$ws.Range("E4").Value = '  -0.01%  '

$ws.Range("D5").Value = '''267.98'
$ws.Range("E5").Value = '  +4.62%  '

$ws.Range("D6").Value = '''85.89'
$ws.Range("E6").Value = '  +11.23%  '

$ws.Range("E7").Value = '  +0.35%  '

$ws.Range("E8").Value = '  -0.02%  '

$ws.Range("D9").Value = '''0.603'
$ws.Range("E9").Value = '  +1.47%  '

$ws.Range("D10").Value = '''45.96'
$ws.Range("E10").Value = '  +8.47%  '

$ws.Range("D11").Value = '''0.0920'
$ws.Range("E11").Value = '  +1.43%  '

$ws.Range("D12").Value = '''7.51'
$ws.Range("E12").Value = '  +6.63%  '

$ws.Range("D13").Value = '''0.104'
$ws.Range("E13").Value = '  +1.29%  '

$ws.Range("D14").Value = '2.546.25'
$ws.Range("E14").Value = '  -0.18%  '

$ws.Range("D15").Value = '''14.59'
$ws.Range("E15").Value = '  +0.88%  '

$ws.Range("D16").Value = '2.209.60'
$ws.Range("E16").Value = '  +0.10%  '

$ws.Range("D17").Value = '''0.785'
$ws.Range("E17").Value = '  +0.25%  '

$ws.Range("D18").Value = '43.680.98'
$ws.Range("E18").Value = '  +2.08%  '

$ws.Range("E19").Value = '  +0.68%  '

$ws.Range("D20").Value = '''5.99'
$ws.Range("E20").Value = '  +0.21%  '

$ws.Range("D21").Value = '''69.87'
$ws.Range("E21").Value = '  -1.69%  '

$ws.Range("D22").Value = '''2.37'
$ws.Range("E22").Value = '  +5.00%  '

$ws.Range("D23").Value = '''231.79'
$ws.Range("E23").Value = '  +0.71%  '

$ws.Range("B24").Value = 'InternetComputer(DFINITY)'
$ws.Range("C24").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D24").Value = '''8.86'
$ws.Range("E24").Value = '  -4.80%  '

$ws.Range("B25").Value = 'PancakeSwap'
$ws.Range("C25").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D25").Value = '''2.66'
$ws.Range("E25").Value = '  +20.87%  '

$ws.Range("E26").Value = '  +0.04%  '

$ws.Range("D27").Value = '''10.78'
$ws.Range("E27").Value = '  +0.14%  '

$ws.Range("E28").Value = '  +5.74%  '

$ws.Range("D29").Value = '''39.10'
$ws.Range("E29").Value = '  -8.45%  '

$ws.Range("D30").Value = '''2.24'
$ws.Range("E30").Value = '  +1.35%  '

$ws.Range("D31").Value = '''175.51'
$ws.Range("E31").Value = '  +1.58%  '

$ws.Range("D32").Value = '''0.0892'
$ws.Range("E32").Value = '  +1.84%  '

$ws.Range("D33").Value = '''20.53'
$ws.Range("E33").Value = '  +0.72%  '

$ws.Range("D34").Value = '''5.40'
$ws.Range("E34").Value = '  +3.28%  '

$ws.Range("E35").Value = '  +2.02%  '

$ws.Range("E36").Value = '  +2.65%  '

$ws.Range("D37").Value = '''0.0355'
$ws.Range("E37").Value = '  -0.94%  '

$ws.Range("D38").Value = '''4.34'
$ws.Range("E38").Value = '  +0.16%  '

$ws.Range("D39").Value = '''3.28'
$ws.Range("E39").Value = '  +16.74%  '

$ws.Range("D40").Value = '''12.35'
$ws.Range("E40").Value = '  -5.98%  '

$ws.Range("D41").Value = '''64.87'
$ws.Range("E41").Value = '  +7.83%  '

$ws.Range("E42").Value = '  -0.87%  '

$ws.Range("D43").Value = '''0.203'
$ws.Range("E43").Value = '  +1.05%  '

$ws.Range("D44").Value = '''5.36'
$ws.Range("E44").Value = '  +1.02%  '

$ws.Range("E45").Value = '  +1.21%  '

$ws.Range("D46").Value = '''8.36'
$ws.Range("E46").Value = '  -0.41%  '

$ws.Range("D47").Value = '''100.18'
$ws.Range("E47").Value = '  -3.05%  '

$ws.Range("E48").Value = '  +6.15%  '

$ws.Range("E49").Value = '  +0.54%  '

$ws.Range("E50").Value = '  -7.01%  '

$ws.Range("D51").Value = '''1.50'
$ws.Range("E51").Value = '  +4.81%  '
